$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells: new "team record" columns
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match header style (bold, centered, bordered) used by the rest of row 1
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Data rows 2-45: team record for every player row
$ws.Range("AD2:AD45").Value = 91
$ws.Range("AE2:AE45").Value = 71
$ws.Range("AF2:AF45").Value = 0
